# "Added Roslyn Model and Test Result" ------------------------------------
# Both sheets ("Label Prediction Results" and "Assignee Prediction Results")
# get the same shape of edit:
#   - H2 label changes from "Epochs(=30)" to "Epochs(=30-60)"
#   - H3 label changes from "Train/Test Data (80/20)" moved up one slot
#     (shared-string index shuffles as a side effect of the old "Epochs(=30)"
#     string being fully replaced/garbage-collected)
#   - row 4 ("Roslyn" / "Title") gains Precision/Recall/F-Measure/Accuracy
#     numbers that were previously blank
#   - row 5 ("Title + Description") gains the same four numbers
# Sheet1 additionally grows two brand-new rows for a second new dataset
# ("corefx"), and its backing table is resized to include them.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Label Prediction Results
$ws2 = $wb.Worksheets.Item(2)   # Assignee Prediction Results

# ---------------------------------------------------------------------
# Sheet 1: Label Prediction Results
# ---------------------------------------------------------------------

# Epochs note increased from 30 to 30-60
$ws1.Range("H2").Value = "Epochs(=30-60)"
$ws1.Range("H3").Value = "Train/Test Data (80/20)"

# Roslyn / Title row now has real numbers
$ws1.Range("A4").Value = "Roslyn"
$ws1.Range("B4").Value = "Title"
$ws1.Range("C4").Value = 3
$ws1.Range("D4").Value = 2
$ws1.Range("E4").Value = 2
$ws1.Range("F4").Value = 6

# Roslyn / Title + Description row now has real numbers
$ws1.Range("B5").Value = "Title + Description"
$ws1.Range("C5").Value = 2
$ws1.Range("D5").Value = 2
$ws1.Range("E5").Value = 2
$ws1.Range("F5").Value = 9

# New dataset "corefx" added as two new rows
$ws1.Range("A6").Value = "corefx"
$ws1.Range("B6").Value = "Title"
$ws1.Range("B7").Value = "Title + Description"

# Grow Table2 (sheet1's table) to cover the two new rows
$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:F7"))

# ---------------------------------------------------------------------
# Sheet 2: Assignee Prediction Results
# ---------------------------------------------------------------------

$ws2.Range("H2").Value = "Epochs(=30-60)"
$ws2.Range("H3").Value = "Train/Test Data (80/20)"

$ws2.Range("A4").Value = "Roslyn"
$ws2.Range("B4").Value = "Title"
$ws2.Range("C4").Value = 1
$ws2.Range("D4").Value = 2
$ws2.Range("E4").Value = 1
$ws2.Range("F4").Value = 6

$ws2.Range("B5").Value = "Title + Description"
$ws2.Range("C5").Value = 4
$ws2.Range("D5").Value = 4
$ws2.Range("E5").Value = 3
$ws2.Range("F5").Value = 12

# The author's cursor ended up on C14 of the second sheet; re-activate
# sheet 1 afterwards so it stays the selected tab, matching the workbook.
$ws2.Range("C14").Select()
$ws1.Activate()
